$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 951.73334
$ws.Range("J17").Value = 951.73334
$ws.Range("L17").Value = 2855.20002
$ws.Range("N17").Value = -3191.20002
$ws.Range("H40").Value = 1538.3125
$ws.Range("I40").Value = 1481.08
$ws.Range("K40").Value = 1481.08
$ws.Range("M40").Value = -1306.08
$ws.Range("H76").Value = 1945
$ws.Range("I76").Value = 1945
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 1945
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -1630
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 1945
$ws.Range("I79").Value = 1945
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 1945
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -853
$ws.Range("N79").ClearContents()
$ws.Range("H116").Value = 7349.875
$ws.Range("I116").Value = 11998.667
$ws.Range("K116").Value = 11998.667
$ws.Range("M116").Value = -8556.666999999999
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H138").Value = 4572.933
$ws.Range("I138").Value = 4559.8
$ws.Range("J138").Value = 4579.5
$ws.Range("K138").Value = 13679.4
$ws.Range("L138").Value = 13738.5
$ws.Range("M138").Value = -8539.400000000001
$ws.Range("N138").Value = -24018.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 20999
$ws.Range("J23").Value = 20999
$ws.Range("L23").Value = 20999
$ws.Range("N23").Value = -21517
$ws.Range("H61").Value = 2360.4614
$ws.Range("I61").Value = 2372
$ws.Range("J61").Value = 2297
$ws.Range("K61").Value = 2372
$ws.Range("L61").Value = 2297
$ws.Range("M61").Value = -2160
$ws.Range("N61").Value = -2721
$ws.Range("H74").Value = 1039.9166
$ws.Range("I74").Value = 1190
$ws.Range("J74").Value = 289.5
$ws.Range("K74").Value = 1190
$ws.Range("L74").Value = 289.5
$ws.Range("M74").Value = -316
$ws.Range("N74").Value = -2037.5
$ws.Range("H77").Value = 1039.9166
$ws.Range("I77").Value = 1190
$ws.Range("J77").Value = 289.5
$ws.Range("K77").Value = 5950
$ws.Range("L77").Value = 1447.5
$ws.Range("M77").Value = -1582
$ws.Range("N77").Value = -10183.5
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
$ws.Range("H110").Value = 2861.3333
$ws.Range("I110").Value = 2861.3333
$ws.Range("K110").Value = 2861.3333
$ws.Range("M110").Value = -816.3332999999998
$ws.Range("H124").Value = 62499.75
$ws.Range("J124").Value = 62499.75
$ws.Range("L124").Value = 62499.75
$ws.Range("N124").Value = -72319.75
$ws.Range("H132").Value = 1959.6
$ws.Range("J132").Value = 1999
$ws.Range("L132").Value = 5997
$ws.Range("N132").Value = -11057
$ws.Range("H135").Value = 525214.5
$ws.Range("J135").Value = 525214.5
$ws.Range("L135").Value = 525214.5
$ws.Range("N135").Value = -535354.5
$ws.Range("H136").Value = 2360.4614
$ws.Range("I136").Value = 2372
$ws.Range("J136").Value = 2297
$ws.Range("K136").Value = 7116
$ws.Range("L136").Value = 6891
$ws.Range("M136").Value = -4566
$ws.Range("N136").Value = -11991
$ws.Range("H138").Value = 488366.28
$ws.Range("J138").Value = 488366.28
$ws.Range("L138").Value = 488366.28
$ws.Range("N138").Value = -498646.28

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 7999.75
$ws.Range("I16").Value = 5999.5
$ws.Range("J16").Value = 10000
$ws.Range("K16").Value = 5999.5
$ws.Range("L16").Value = 10000
$ws.Range("M16").Value = -5712.5
$ws.Range("N16").Value = -10574
$ws.Range("H22").Value = 196.77777
$ws.Range("I22").Value = 212.83333
$ws.Range("J22").Value = 164.66667
$ws.Range("K22").Value = 212.83333
$ws.Range("L22").Value = 164.66667
$ws.Range("M22").Value = 137.16667
$ws.Range("N22").Value = -864.6666700000001
$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()
$ws.Range("H86").Value = 11062.857
$ws.Range("I86").Value = 11909
$ws.Range("K86").Value = 11909
$ws.Range("M86").Value = -10786
$ws.Range("H89").Value = 11062.857
$ws.Range("I89").Value = 11909
$ws.Range("K89").Value = 59545
$ws.Range("M89").Value = -53929
$ws.Range("H113").Value = 7999.75
$ws.Range("I113").Value = 5999.5
$ws.Range("J113").Value = 10000
$ws.Range("K113").Value = 5999.5
$ws.Range("L113").Value = 10000
$ws.Range("M113").Value = -3829.5
$ws.Range("N113").Value = -14340
$ws.Range("H132").Value = 3584.2856
$ws.Range("I132").Value = 3681.6667
$ws.Range("K132").Value = 11045.0001
$ws.Range("M132").Value = -8515.000100000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 4166.3335
$ws.Range("J132").Value = 4250
$ws.Range("L132").Value = 38250
$ws.Range("N132").Value = -43310
$ws.Range("H134").Value = 1332
$ws.Range("I134").Value = 1332
$ws.Range("K134").Value = 3996
$ws.Range("M134").Value = 1074

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2486.5833
$ws.Range("I102").Value = 2539.6
$ws.Range("J102").Value = 2221.5
$ws.Range("K102").Value = 2539.6
$ws.Range("L102").Value = 2221.5
$ws.Range("M102").Value = -917.5999999999999
$ws.Range("N102").Value = -5465.5
$ws.Range("H113").Value = 897.9
$ws.Range("I113").Value = 857
$ws.Range("J113").Value = 938.8
$ws.Range("K113").Value = 857
$ws.Range("L113").Value = 938.8
$ws.Range("M113").Value = 1313
$ws.Range("N113").Value = -5278.8
$ws.Range("H126").Value = 1995.7142
$ws.Range("I126").Value = 1995.3334
$ws.Range("K126").Value = 5986.0002
$ws.Range("M126").Value = -3516.0002
$ws.Range("H132").Value = 8522.125
$ws.Range("I132").Value = 8522.125
$ws.Range("K132").Value = 25566.375
$ws.Range("M132").Value = -23036.375
$ws.Range("H140").Value = 143022.5
$ws.Range("J140").Value = 143022.5
$ws.Range("L140").Value = 143022.5
$ws.Range("N140").Value = -153382.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5999.6665
$ws.Range("I40").Value = 5999
$ws.Range("J40").Value = 6000
$ws.Range("K40").Value = 5999
$ws.Range("L40").Value = 6000
$ws.Range("M40").Value = -5863
$ws.Range("N40").Value = -6272
$ws.Range("H61").Value = 1698.8
$ws.Range("I61").Value = 1499
$ws.Range("J61").Value = 1998.5
$ws.Range("K61").Value = 1499
$ws.Range("L61").Value = 1998.5
$ws.Range("M61").Value = -1297
$ws.Range("N61").Value = -2402.5
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("M75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("M78").ClearContents()
$ws.Range("H113").Value = 1698.8
$ws.Range("I113").Value = 1499
$ws.Range("J113").Value = 1998.5
$ws.Range("K113").Value = 1499
$ws.Range("L113").Value = 1998.5
$ws.Range("M113").Value = 671
$ws.Range("N113").Value = -6338.5
$ws.Range("H122").Value = 5694.567
$ws.Range("I122").Value = 4655.9414
$ws.Range("J122").Value = 7052.769
$ws.Range("K122").Value = 13967.8242
$ws.Range("L122").Value = 21158.307
$ws.Range("M122").Value = -11517.8242
$ws.Range("N122").Value = -26058.307
$ws.Range("H127").Value = 49999
$ws.Range("J127").Value = 49999
$ws.Range("L127").Value = 49999
$ws.Range("N127").Value = -59919
$ws.Range("H132").Value = 12956.6
$ws.Range("I132").Value = 18327.666
$ws.Range("K132").Value = 54982.99800000001
$ws.Range("M132").Value = -52452.99800000001
$ws.Range("H136").Value = 3934.2856
$ws.Range("I136").Value = 3934.2856
$ws.Range("K136").Value = 11802.8568
$ws.Range("M136").Value = -9252.856800000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 19643.428
$ws.Range("I41").Value = 19663.5
$ws.Range("K41").Value = 19663.5
$ws.Range("M41").Value = -19273.5

